$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 78995
$ws.Range("J3").Value = 78995
$ws.Range("L3").Value = 78995
$ws.Range("N3").Value = -79223
$ws.Range("H55").Value = 110.14286
$ws.Range("I55").Value = 100.25
$ws.Range("K55").Value = 100.25
$ws.Range("M55").Value = 113.75
$ws.Range("H57").Value = 64372.5
$ws.Range("J57").Value = 64372.5
$ws.Range("L57").Value = 193117.5
$ws.Range("N57").Value = -194115.5
$ws.Range("H62").Value = 6938.4
$ws.Range("I62").Value = 5861
$ws.Range("K62").Value = 5861
$ws.Range("M62").Value = -5237
$ws.Range("H65").Value = 6938.4
$ws.Range("I65").Value = 5861
$ws.Range("K65").Value = 29305
$ws.Range("M65").Value = -26185
$ws.Range("H87").Value = 108085.836
$ws.Range("J87").Value = 108085.836
$ws.Range("L87").Value = 108085.836
$ws.Range("N87").Value = -110581.836
$ws.Range("H90").Value = 108085.836
$ws.Range("J90").Value = 108085.836
$ws.Range("L90").Value = 324257.508
$ws.Range("N90").Value = -336737.508
$ws.Range("H101").Value = 1910.5
$ws.Range("I101").Value = 2022.8182
$ws.Range("K101").Value = 6068.4546
$ws.Range("M101").Value = -4446.4546
$ws.Range("H102").Value = 78995
$ws.Range("J102").Value = 78995
$ws.Range("L102").Value = 78995
$ws.Range("N102").Value = -85485
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("H105").Value = 58329.668
$ws.Range("J105").Value = 58329.668
$ws.Range("L105").Value = 58329.668
$ws.Range("N105").Value = -65317.668
$ws.Range("H110").Value = 148000
$ws.Range("J110").Value = 148000
$ws.Range("L110").Value = 148000
$ws.Range("N110").Value = -156180
$ws.Range("H114").Value = 69995
$ws.Range("J114").Value = 69995
$ws.Range("L114").Value = 69995
$ws.Range("N114").Value = -78673
$ws.Range("H115").Value = 894
$ws.Range("I115").Value = 842.5
$ws.Range("J115").Value = 1100
$ws.Range("K115").Value = 2527.5
$ws.Range("L115").Value = 3300
$ws.Range("M115").Value = -960.5
$ws.Range("N115").Value = -6434
$ws.Range("H117").Value = 164000
$ws.Range("J117").Value = 164000
$ws.Range("L117").Value = 164000
$ws.Range("N117").Value = -173178

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1304.4667
$ws.Range("I2").Value = 1131.0834
$ws.Range("J2").Value = 1998
$ws.Range("K2").Value = 1131.0834
$ws.Range("L2").Value = 1998
$ws.Range("M2").Value = -1018.0834
$ws.Range("N2").Value = -2224
$ws.Range("H32").Value = 5451.41
$ws.Range("I32").Value = 2389.2876
$ws.Range("J32").Value = 17699.9
$ws.Range("K32").Value = 2389.2876
$ws.Range("L32").Value = 17699.9
$ws.Range("M32").Value = -2102.2876
$ws.Range("N32").Value = -18273.9
$ws.Range("H61").Value = 30004378
$ws.Range("I61").Value = 23812510
$ws.Range("J61").Value = 62511680
$ws.Range("K61").Value = 23812510
$ws.Range("L61").Value = 62511680
$ws.Range("M61").Value = -23812298
$ws.Range("N61").Value = -62512104
$ws.Range("H101").Value = 54998.43
$ws.Range("J101").Value = 54998.43
$ws.Range("L101").Value = 54998.43
$ws.Range("N101").Value = -61488.43
$ws.Range("H102").Value = 3281.48
$ws.Range("I102").Value = 3809.524
$ws.Range("K102").Value = 3809.524
$ws.Range("M102").Value = -2187.524
$ws.Range("H116").Value = 1304.4667
$ws.Range("I116").Value = 1131.0834
$ws.Range("J116").Value = 1998
$ws.Range("K116").Value = 1131.0834
$ws.Range("L116").Value = 1998
$ws.Range("M116").Value = 1162.9166
$ws.Range("N116").Value = -6586
$ws.Range("H136").Value = 30004378
$ws.Range("I136").Value = 23812510
$ws.Range("J136").Value = 62511680
$ws.Range("K136").Value = 71437530
$ws.Range("L136").Value = 187535040
$ws.Range("M136").Value = -71434980
$ws.Range("N136").Value = -187540140

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1304.4667
$ws.Range("I3").Value = 1131.0834
$ws.Range("J3").Value = 1998
$ws.Range("K3").Value = 1131.0834
$ws.Range("L3").Value = 1998
$ws.Range("M3").Value = -1017.0834
$ws.Range("N3").Value = -2226
$ws.Range("H94").Value = 1514.9615
$ws.Range("I94").Value = 612.875
$ws.Range("K94").Value = 612.875
$ws.Range("M94").Value = -161.875
$ws.Range("H105").Value = 2168.8518
$ws.Range("I105").Value = 1573.2858
$ws.Range("K105").Value = 1573.2858
$ws.Range("M105").Value = 173.7141999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 610491.5600000001
$ws.Range("I31").Value = 11378.533
$ws.Range("K31").Value = 11378.533
$ws.Range("M31").Value = -11083.533
$ws.Range("H34").Value = 610491.5600000001
$ws.Range("I34").Value = 11378.533
$ws.Range("K34").Value = 11378.533
$ws.Range("M34").Value = -11176.533
$ws.Range("H57").Value = 25000
$ws.Range("J57").Value = 25000
$ws.Range("L57").Value = 25000
$ws.Range("N57").Value = -26120
$ws.Range("H58").Value = 5166.769
$ws.Range("I58").Value = 3583.25
$ws.Range("J58").Value = 7700.4
$ws.Range("K58").Value = 3583.25
$ws.Range("L58").Value = 7700.4
$ws.Range("M58").Value = -3380.25
$ws.Range("N58").Value = -8106.4
$ws.Range("H99").Value = 2851.7778
$ws.Range("I99").Value = 2543.4
$ws.Range("K99").Value = 2543.4
$ws.Range("M99").Value = -1045.4
$ws.Range("H122").Value = 2955.75
$ws.Range("I122").Value = 2941
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8823
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6373
$ws.Range("N122").Value = -13900
$ws.Range("H126").Value = 2851.7778
$ws.Range("I126").Value = 2543.4
$ws.Range("K126").Value = 7630.200000000001
$ws.Range("M126").Value = -5160.200000000001
$ws.Range("H136").Value = 5166.769
$ws.Range("I136").Value = 3583.25
$ws.Range("J136").Value = 7700.4
$ws.Range("K136").Value = 10749.75
$ws.Range("L136").Value = 23101.2
$ws.Range("M136").Value = -8199.75
$ws.Range("N136").Value = -28201.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7722464.5
$ws.Range("I4").Value = 6030912.5
$ws.Range("K4").Value = 18092737.5
$ws.Range("M4").Value = -18092625.5
$ws.Range("H7").Value = 1412.2858
$ws.Range("I7").Value = 376.6
$ws.Range("K7").Value = 1129.8
$ws.Range("M7").Value = -1017.8
$ws.Range("H37").Value = 84499
$ws.Range("J37").Value = 84499
$ws.Range("L37").Value = 253497
$ws.Range("N37").Value = -253721
$ws.Range("H38").Value = 346.44446
$ws.Range("I38").Value = 569.75
$ws.Range("J38").Value = 167.8
$ws.Range("K38").Value = 1709.25
$ws.Range("L38").Value = 503.4
$ws.Range("M38").Value = -1362.25
$ws.Range("N38").Value = -1197.4
$ws.Range("H134").Value = 4336.125
$ws.Range("I134").Value = 1448.5
$ws.Range("K134").Value = 4345.5
$ws.Range("M134").Value = 724.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 56998.332
$ws.Range("J95").Value = 56998.332
$ws.Range("L95").Value = 56998.332
$ws.Range("N95").Value = -62490.332
$ws.Range("H110").Value = 101056.664
$ws.Range("J110").Value = 101056.664
$ws.Range("L110").Value = 101056.664
$ws.Range("N110").Value = -109236.664
$ws.Range("H122").Value = 3597.4614
$ws.Range("I122").Value = 2659.3333
$ws.Range("K122").Value = 7977.999899999999
$ws.Range("M122").Value = -5527.999899999999
$ws.Range("H124").Value = 192973
$ws.Range("J124").Value = 192973
$ws.Range("L124").Value = 192973
$ws.Range("N124").Value = -202793
$ws.Range("H126").Value = 3329.5264
$ws.Range("I126").Value = 2481.6155
$ws.Range("K126").Value = 7444.8465
$ws.Range("M126").Value = -4974.8465
$ws.Range("H132").Value = 28574892
$ws.Range("I132").Value = 31253382
$ws.Range("K132").Value = 93760146
$ws.Range("M132").Value = -93757616

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3844.75
$ws.Range("I46").Value = 3231.7
$ws.Range("J46").Value = 4866.5
$ws.Range("K46").Value = 3231.7
$ws.Range("L46").Value = 4866.5
$ws.Range("M46").Value = -3043.7
$ws.Range("N46").Value = -5242.5
$ws.Range("H132").Value = 322508.72
$ws.Range("I132").Value = 11399.538
$ws.Range("J132").Value = 1670648.5
$ws.Range("K132").Value = 34198.614
$ws.Range("L132").Value = 5011945.5
$ws.Range("M132").Value = -31668.614
$ws.Range("N132").Value = -5017005.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 13634
$ws.Range("I52").Value = 13634
$ws.Range("K52").Value = 13634
$ws.Range("M52").Value = -13408
$ws.Range("H126").Value = 2269.5454
$ws.Range("I126").Value = 2296.5
$ws.Range("K126").Value = 6889.5
$ws.Range("M126").Value = -4419.5
